# Auto update: 2025-12-03 03:05:24
# Refreshes the 미장 반도체 분석 sheet for the 2025-12-03 run:
#  - advances the report date (A2:A6) from 2025-12-01 to 2025-12-03
#  - QUALCOMM now lands in row 4 (was AMD) and AMD now lands in row 6 (was QUALCOMM)
#  - recomputes the per-ticker indicator columns (D..K, N)
#  - rewords the MACRO_SIGNAL banner text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A : date refresh (rows 2-6) ---
# Force text storage (matches the source file, where the date is a literal
# string, not a real date serial) by flipping the format to Text before the
# write, then handing the range back its default style so the cells keep
# the workbook's normal (unstyled) look.
$dateRng = $ws.Range("A2:A6")
$dateRng.NumberFormat = "@"
$dateRng.Value = "2025-12-03"
$dateRng.Style = "Normal"

# --- Row 4 : now QUALCOMM Incorporated / QCOM (was AMD / AMD) ---
$ws.Range("B4").Value = "QUALCOMM Incorporated"
$ws.Range("C4").Value = "QCOM"

# --- Row 6 : now Advanced Micro Devices, Inc. / AMD (was QUALCOMM / QCOM) ---
$ws.Range("B6").Value = "Advanced Micro Devices, Inc."
$ws.Range("C6").Value = "AMD"

# --- Row 2 : Taiwan Semiconductor (TSM) indicator refresh ---
$ws.Range("D2").Value = 292.03
$ws.Range("E2").Value = 50.8
$ws.Range("F2").Value = 2.6
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 56
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 62.8
$ws.Range("N2").Value = 65.32892478746797

# --- Row 3 : ASML Holding (ASML) indicator refresh ---
$ws.Range("D3").Value = 1098.26
$ws.Range("E3").Value = 62.1
$ws.Range("F3").Value = 11.18
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 53
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 55.8
$ws.Range("N3").Value = 65.32892478746797

# --- Row 4 : QUALCOMM (QCOM) indicator refresh ---
$ws.Range("D4").Value = 169.92
$ws.Range("E4").Value = 44.3
$ws.Range("F4").Value = 2.94
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 55.6
$ws.Range("N4").Value = 65.32892478746797

# --- Row 5 : NVIDIA Corporation (NVDA) indicator refresh ---
$ws.Range("D5").Value = 181.43
$ws.Range("E5").Value = 38.5
$ws.Range("F5").Value = -0.61
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 66
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 49
$ws.Range("N5").Value = 65.32892478746797

# --- Row 6 : AMD indicator refresh ---
$ws.Range("D6").Value = 217.33
$ws.Range("E6").Value = 41
$ws.Range("F6").Value = 1.06
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 46
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 56
$ws.Range("K6").Value = 44.6
$ws.Range("N6").Value = 65.32892478746797

# --- MACRO_SIGNAL text update (O column shares this string across all rows) ---
$ws.Range("O2:O6").Value = "🟢 상승 우위 (다소 완화)"
